# Updated symbol list on Fri Dec 30 04:00:45 UTC 2022 with GitHub Actions
#
# Refreshes the scraped "Price" (D), "Volume(1h)" (E) and "Hora" (G) columns
# on the crypto-ranking sheet with the latest pull. Price/Hora values are
# numeric-looking text (scraper writes them as literal strings, preserving
# trailing zeros / leading zeros, e.g. "0.001520"), so a leading apostrophe
# forces Excel to keep them as text instead of silently coercing to a
# number (which would normalize "0.001520" -> 0.00152 and "4" -> 4, losing
# the original formatting). The cell's pre-existing Style is restored right
# after so the quote-prefix doesn't leave a stray numeric-as-text marker on
# the cell's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Ref = "D2"; Value = "245.34"; Numeric = $true },
    @{ Ref = "G2"; Value = "4"; Numeric = $true },
    @{ Ref = "D3"; Value = "24.14"; Numeric = $true },
    @{ Ref = "G3"; Value = "4"; Numeric = $true },
    @{ Ref = "D4"; Value = "5.251"; Numeric = $true },
    @{ Ref = "G4"; Value = "4"; Numeric = $true },
    @{ Ref = "D5"; Value = "0.05791"; Numeric = $true },
    @{ Ref = "G5"; Value = "4"; Numeric = $true },
    @{ Ref = "G6"; Value = "4"; Numeric = $true },
    @{ Ref = "D7"; Value = "3.133"; Numeric = $true },
    @{ Ref = "G7"; Value = "4"; Numeric = $true },
    @{ Ref = "D8"; Value = "0.8181"; Numeric = $true },
    @{ Ref = "G8"; Value = "4"; Numeric = $true },
    @{ Ref = "D9"; Value = "0.8584"; Numeric = $true },
    @{ Ref = "G9"; Value = "4"; Numeric = $true },
    @{ Ref = "G10"; Value = "4"; Numeric = $true },
    @{ Ref = "D11"; Value = "0.06943"; Numeric = $true },
    @{ Ref = "G11"; Value = "4"; Numeric = $true },
    @{ Ref = "D12"; Value = "0.03207"; Numeric = $true },
    @{ Ref = "G12"; Value = "4"; Numeric = $true },
    @{ Ref = "D13"; Value = "0.02877"; Numeric = $true },
    @{ Ref = "G13"; Value = "4"; Numeric = $true },
    @{ Ref = "D14"; Value = "0.09384"; Numeric = $true },
    @{ Ref = "G14"; Value = "4"; Numeric = $true },
    @{ Ref = "D15"; Value = "3.749"; Numeric = $true },
    @{ Ref = "G15"; Value = "4"; Numeric = $true },
    @{ Ref = "D16"; Value = "0.001520"; Numeric = $true },
    @{ Ref = "G16"; Value = "4"; Numeric = $true },
    @{ Ref = "D17"; Value = "0.04706"; Numeric = $true },
    @{ Ref = "G17"; Value = "4"; Numeric = $true },
    @{ Ref = "D18"; Value = "0.0005961"; Numeric = $true },
    @{ Ref = "E18"; Value = "17OneONE"; Numeric = $false },
    @{ Ref = "G18"; Value = "4"; Numeric = $true },
    @{ Ref = "D19"; Value = "0.006279"; Numeric = $true },
    @{ Ref = "G19"; Value = "4"; Numeric = $true },
    @{ Ref = "D20"; Value = "0.001235"; Numeric = $true },
    @{ Ref = "G20"; Value = "4"; Numeric = $true },
    @{ Ref = "D21"; Value = "0.004611"; Numeric = $true },
    @{ Ref = "G21"; Value = "4"; Numeric = $true },
    @{ Ref = "E22"; Value = "21NitroExNTXWorstin24h"; Numeric = $false },
    @{ Ref = "G22"; Value = "4"; Numeric = $true },
    @{ Ref = "D23"; Value = "3.508"; Numeric = $true },
    @{ Ref = "G23"; Value = "4"; Numeric = $true },
    @{ Ref = "D24"; Value = "2.121"; Numeric = $true },
    @{ Ref = "G24"; Value = "4"; Numeric = $true },
    @{ Ref = "D25"; Value = "0.3191"; Numeric = $true },
    @{ Ref = "G25"; Value = "4"; Numeric = $true },
    @{ Ref = "G26"; Value = "4"; Numeric = $true },
    @{ Ref = "G27"; Value = "4"; Numeric = $true },
    @{ Ref = "D28"; Value = "0.0002329"; Numeric = $true },
    @{ Ref = "G28"; Value = "4"; Numeric = $true },
    @{ Ref = "G29"; Value = "4"; Numeric = $true },
    @{ Ref = "G30"; Value = "4"; Numeric = $true },
    @{ Ref = "G31"; Value = "4"; Numeric = $true },
    @{ Ref = "G32"; Value = "4"; Numeric = $true },
    @{ Ref = "G33"; Value = "4"; Numeric = $true },
    @{ Ref = "G34"; Value = "4"; Numeric = $true },
    @{ Ref = "G35"; Value = "4"; Numeric = $true },
    @{ Ref = "G36"; Value = "4"; Numeric = $true },
    @{ Ref = "G37"; Value = "4"; Numeric = $true },
    @{ Ref = "G38"; Value = "4"; Numeric = $true },
    @{ Ref = "G39"; Value = "4"; Numeric = $true },
    @{ Ref = "D40"; Value = "0.03652"; Numeric = $true },
    @{ Ref = "G40"; Value = "4"; Numeric = $true },
    @{ Ref = "D41"; Value = "0.006281"; Numeric = $true },
    @{ Ref = "E41"; Value = "40KickTokenKICKBestin24h"; Numeric = $false },
    @{ Ref = "G41"; Value = "4"; Numeric = $true },
    @{ Ref = "G42"; Value = "4"; Numeric = $true },
    @{ Ref = "E43"; Value = "42CEJICEJI"; Numeric = $false },
    @{ Ref = "G43"; Value = "4"; Numeric = $true },
    @{ Ref = "D44"; Value = "0.008058"; Numeric = $true },
    @{ Ref = "G44"; Value = "4"; Numeric = $true },
    @{ Ref = "G45"; Value = "4"; Numeric = $true },
    @{ Ref = "G46"; Value = "4"; Numeric = $true },
    @{ Ref = "G47"; Value = "4"; Numeric = $true },
    @{ Ref = "D48"; Value = "0.002339"; Numeric = $true },
    @{ Ref = "G48"; Value = "4"; Numeric = $true },
    @{ Ref = "G49"; Value = "4"; Numeric = $true },
    @{ Ref = "G50"; Value = "4"; Numeric = $true },
    @{ Ref = "G51"; Value = "4"; Numeric = $true }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Ref)
    if ($change.Numeric) {
        $origStyle = $cell.Style
        $cell.Value = "'" + $change.Value
        $cell.Style = $origStyle
    } else {
        $cell.Value = $change.Value
    }
}
